$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.406915068626404
$ws.Range("B1").Value = 1.48189389705658
$ws.Range("C1").Value = 1.3862384557724
$ws.Range("D1").Value = 1.448469161987305
$ws.Range("E1").Value = 1.049007534980774
